# Insert a new data row at row 48 (shifts existing rows 48:102 down to 49:103)
# and populate it with the new weekly record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(48).Insert()

$ws.Cells.Item(48, 1).Value = 11
$ws.Cells.Item(48, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(48, 3).Value = "Bíobío"
$ws.Cells.Item(48, 4).Value = 44763
$ws.Cells.Item(48, 5).Value = 8
$ws.Cells.Item(48, 6).Value = 100112001
$ws.Cells.Item(48, 7).Value = "Berenjena"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 180
$ws.Cells.Item(48, 11).Value = 11000
$ws.Cells.Item(48, 12).Value = 12000
$ws.Cells.Item(48, 13).Value = 11444
$ws.Cells.Item(48, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(48, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(48, 16).Value = 191
$ws.Cells.Item(48, 17).Value = 60
$ws.Cells.Item(48, 18).Value = "Hortaliza"
